$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "user4" row (row 6) entirely - the test fixture was trimmed
# down to 3 users (admin, user1, user2, user3) for the auth/CRUD robot
# tests, so the last data row goes away and everything below it shifts up.
$ws.Rows("6:6").Delete() | Out-Null

# The password column no longer stores a bcrypt hash for the remaining
# rows - it now just mirrors the username column.
$ws.Range("B2").Value = "admin"
$ws.Range("B3").Value = "user1"
$ws.Range("B4").Value = "user2"
$ws.Range("B5").Value = "user3"

# Leave the selection where the author left it when they saved the file.
$ws.Range("B9").Select() | Out-Null
